# --implement getRowContains and test it
#
# Replaces the two ad-hoc probe cells (F6, F8) and the stray E2 cell with two
# proper new data rows (3 and 4) that exercise RunTestCase_ICSLogin /
# RunTestCase_OnlineStoreLogin, each with its own mailto hyperlink on column D
# (matching the existing D2 pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old scratch cells that are no longer needed.
$ws.Range("E2").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F8").ClearContents()

# Row 3: D / RunTestCase_ICSLogin / user|password / admin|Admin
$ws.Range("A3").Value = "D"
$ws.Range("B3").Value = "RunTestCase_ICSLogin"
$ws.Range("C3").Value = "user|password"
$ws.Range("D3").Value = "admin|Admin"

# Row 4: R / RunTestCase_OnlineStoreLogin / user|password / testuser_1234|Test@123
$ws.Range("A4").Value = "R"
$ws.Range("B4").Value = "RunTestCase_OnlineStoreLogin"
$ws.Range("C4").Value = "user|password"
$ws.Range("D4").Value = "testuser_1234|Test@123"

# New row's credential cell gets a mailto hyperlink, just like D2.
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:testuser_1234%7CTest@123")

# Excel auto-applies its built-in "Hyperlink" cell style when a hyperlink is
# added; the source workbook's existing D2 hyperlink cell carries no such
# style, so reset D4 back to Normal and drop the now-unused style def.
$ws.Range("D4").Style = "Normal"
$wb.Styles.Item("Hyperlink").Delete()

[void]$ws.Range("B4").Select()
